# Update LR-pairs data for B2m-Gm11127 with new TPM-derived values.
# Adds MuSCs as a "Sending cluster" (rows 14-17) completing the full 4x4
# (ECs, FAPs, MuSCs, Resolving-Mac) x (ECs, FAPs, MuSCs, Resolving-Mac) grid
# and refreshes all numeric columns (E:T) for every existing row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "B2m"
$ws.Range("C2").Value2 = "Gm11127"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1809.511841
$ws.Range("H2").Value2 = 5428.535523
$ws.Range("I2").Value2 = 0.2989288062794902
$ws.Range("J2").Value2 = 0.2989288062794903
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 0.4271833333333333
$ws.Range("N2").Value2 = 1.28155
$ws.Range("O2").Value2 = 0.7619379997514806
$ws.Range("P2").Value2 = 0.7619379997514806
$ws.Range("Q2").Value2 = 772.9932999445167
$ws.Range("R2").Value2 = 6956.93969950065
$ws.Range("S2").Value2 = 0.2277652167246926
$ws.Range("T2").Value2 = 0.2277652167246926

# Row 3: ECs -> FAPs
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "B2m"
$ws.Range("C3").Value2 = "Gm11127"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1809.511841
$ws.Range("H3").Value2 = 5428.535523
$ws.Range("I3").Value2 = 0.2989288062794902
$ws.Range("J3").Value2 = 0.2989288062794903
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.03181266666666666
$ws.Range("N3").Value2 = 0.095438
$ws.Range("O3").Value2 = 0.05674210044109227
$ws.Range("P3").Value2 = 0.05674210044109228
$ws.Range("Q3").Value2 = 57.56539702711933
$ws.Range("R3").Value2 = 518.088573244074
$ws.Range("S3").Value2 = 0.01696184835064665
$ws.Range("T3").Value2 = 0.01696184835064665

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "B2m"
$ws.Range("C4").Value2 = "Gm11127"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1809.511841
$ws.Range("H4").Value2 = 5428.535523
$ws.Range("I4").Value2 = 0.2989288062794902
$ws.Range("J4").Value2 = 0.2989288062794903
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.02317166666666666
$ws.Range("N4").Value2 = 0.069515
$ws.Range("O4").Value2 = 0.04132973356694953
$ws.Range("P4").Value2 = 0.04132973356694953
$ws.Range("Q4").Value2 = 41.92940520903834
$ws.Range("R4").Value2 = 377.364646881345
$ws.Range("S4").Value2 = 0.0123546479190176
$ws.Range("T4").Value2 = 0.0123546479190176

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value2 = "ECs"
$ws.Range("B5").Value2 = "B2m"
$ws.Range("C5").Value2 = "Gm11127"
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 1809.511841
$ws.Range("H5").Value2 = 5428.535523
$ws.Range("I5").Value2 = 0.2989288062794902
$ws.Range("J5").Value2 = 0.2989288062794903
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.078486
$ws.Range("N5").Value2 = 0.235458
$ws.Range("O5").Value2 = 0.1399901662404776
$ws.Range("P5").Value2 = 0.1399901662404776
$ws.Range("Q5").Value2 = 142.021346352726
$ws.Range("R5").Value2 = 1278.192117174534
$ws.Range("S5").Value2 = 0.04184709328513337
$ws.Range("T5").Value2 = 0.04184709328513338

# Row 6: FAPs -> ECs
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "B2m"
$ws.Range("C6").Value2 = "Gm11127"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 1149.586873333333
$ws.Range("H6").Value2 = 3448.76062
$ws.Range("I6").Value2 = 0.1899101315469672
$ws.Range("J6").Value2 = 0.1899101315469672
$ws.Range("K6").Value2 = 2
$ws.Range("L6").Value2 = 0.6666666666666666
$ws.Range("M6").Value2 = 0.4271833333333333
$ws.Range("N6").Value2 = 1.28155
$ws.Range("O6").Value2 = 0.7619379997514806
$ws.Range("P6").Value2 = 0.7619379997514806
$ws.Range("Q6").Value2 = 491.0843525067778
$ws.Range("R6").Value2 = 4419.759172561
$ws.Range("S6").Value2 = 0.1446997457634367
$ws.Range("T6").Value2 = 0.1446997457634367

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "B2m"
$ws.Range("C7").Value2 = "Gm11127"
$ws.Range("D7").Value2 = "FAPs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 1149.586873333333
$ws.Range("H7").Value2 = 3448.76062
$ws.Range("I7").Value2 = 0.1899101315469672
$ws.Range("J7").Value2 = 0.1899101315469672
$ws.Range("K7").Value2 = 1
$ws.Range("L7").Value2 = 0.3333333333333333
$ws.Range("M7").Value2 = 0.03181266666666666
$ws.Range("N7").Value2 = 0.095438
$ws.Range("O7").Value2 = 0.05674210044109227
$ws.Range("P7").Value2 = 0.05674210044109228
$ws.Range("Q7").Value2 = 36.57142400572889
$ws.Range("R7").Value2 = 329.14281605156
$ws.Range("S7").Value2 = 0.01077589975901906
$ws.Range("T7").Value2 = 0.01077589975901906

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value2 = "FAPs"
$ws.Range("B8").Value2 = "B2m"
$ws.Range("C8").Value2 = "Gm11127"
$ws.Range("D8").Value2 = "MuSCs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 1149.586873333333
$ws.Range("H8").Value2 = 3448.76062
$ws.Range("I8").Value2 = 0.1899101315469672
$ws.Range("J8").Value2 = 0.1899101315469672
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.02317166666666666
$ws.Range("N8").Value2 = 0.069515
$ws.Range("O8").Value2 = 0.04132973356694953
$ws.Range("P8").Value2 = 0.04132973356694953
$ws.Range("Q8").Value2 = 26.63784383325556
$ws.Range("R8").Value2 = 239.7405944993
$ws.Range("S8").Value2 = 0.007848935138500492
$ws.Range("T8").Value2 = 0.007848935138500492

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value2 = "FAPs"
$ws.Range("B9").Value2 = "B2m"
$ws.Range("C9").Value2 = "Gm11127"
$ws.Range("D9").Value2 = "Resolving-Mac"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 1149.586873333333
$ws.Range("H9").Value2 = 3448.76062
$ws.Range("I9").Value2 = 0.1899101315469672
$ws.Range("J9").Value2 = 0.1899101315469672
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.078486
$ws.Range("N9").Value2 = 0.235458
$ws.Range("O9").Value2 = 0.1399901662404776
$ws.Range("P9").Value2 = 0.1399901662404776
$ws.Range("Q9").Value2 = 90.22647534044
$ws.Range("R9").Value2 = 812.03827806396
$ws.Range("S9").Value2 = 0.02658555088601092
$ws.Range("T9").Value2 = 0.02658555088601092

# Row 10: MuSCs -> ECs
$ws.Range("A10").Value2 = "MuSCs"
$ws.Range("B10").Value2 = "B2m"
$ws.Range("C10").Value2 = "Gm11127"
$ws.Range("D10").Value2 = "ECs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 241.4316866666667
$ws.Range("H10").Value2 = 724.29506
$ws.Range("I10").Value2 = 0.03988417442652731
$ws.Range("J10").Value2 = 0.03988417442652733
$ws.Range("K10").Value2 = 2
$ws.Range("L10").Value2 = 0.6666666666666666
$ws.Range("M10").Value2 = 0.4271833333333333
$ws.Range("N10").Value2 = 1.28155
$ws.Range("O10").Value2 = 0.7619379997514806
$ws.Range("P10").Value2 = 0.7619379997514806
$ws.Range("Q10").Value2 = 103.1355926825556
$ws.Range("R10").Value2 = 928.2203341430001
$ws.Range("S10").Value2 = 0.03038926808428738
$ws.Range("T10").Value2 = 0.03038926808428739

# Row 11: MuSCs -> FAPs
$ws.Range("A11").Value2 = "MuSCs"
$ws.Range("B11").Value2 = "B2m"
$ws.Range("C11").Value2 = "Gm11127"
$ws.Range("D11").Value2 = "FAPs"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 241.4316866666667
$ws.Range("H11").Value2 = 724.29506
$ws.Range("I11").Value2 = 0.03988417442652731
$ws.Range("J11").Value2 = 0.03988417442652733
$ws.Range("K11").Value2 = 1
$ws.Range("L11").Value2 = 0.3333333333333333
$ws.Range("M11").Value2 = 0.03181266666666666
$ws.Range("N11").Value2 = 0.095438
$ws.Range("O11").Value2 = 0.05674210044109227
$ws.Range("P11").Value2 = 0.05674210044109228
$ws.Range("Q11").Value2 = 7.680585770697777
$ws.Range("R11").Value2 = 69.12527193628
$ws.Range("S11").Value2 = 0.002263111831320057
$ws.Range("T11").Value2 = 0.002263111831320058

# Row 12: MuSCs -> MuSCs
$ws.Range("A12").Value2 = "MuSCs"
$ws.Range("B12").Value2 = "B2m"
$ws.Range("C12").Value2 = "Gm11127"
$ws.Range("D12").Value2 = "MuSCs"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 241.4316866666667
$ws.Range("H12").Value2 = 724.29506
$ws.Range("I12").Value2 = 0.03988417442652731
$ws.Range("J12").Value2 = 0.03988417442652733
$ws.Range("K12").Value2 = 1
$ws.Range("L12").Value2 = 0.3333333333333333
$ws.Range("M12").Value2 = 0.02317166666666666
$ws.Range("N12").Value2 = 0.069515
$ws.Range("O12").Value2 = 0.04132973356694953
$ws.Range("P12").Value2 = 0.04132973356694953
$ws.Range("Q12").Value2 = 5.59437456621111
$ws.Range("R12").Value2 = 50.3493710959
$ws.Range("S12").Value2 = 0.001648402302586116
$ws.Range("T12").Value2 = 0.001648402302586117

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value2 = "MuSCs"
$ws.Range("B13").Value2 = "B2m"
$ws.Range("C13").Value2 = "Gm11127"
$ws.Range("D13").Value2 = "Resolving-Mac"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 241.4316866666667
$ws.Range("H13").Value2 = 724.29506
$ws.Range("I13").Value2 = 0.03988417442652731
$ws.Range("J13").Value2 = 0.03988417442652733
$ws.Range("K13").Value2 = 1
$ws.Range("L13").Value2 = 0.3333333333333333
$ws.Range("M13").Value2 = 0.078486
$ws.Range("N13").Value2 = 0.235458
$ws.Range("O13").Value2 = 0.1399901662404776
$ws.Range("P13").Value2 = 0.1399901662404776
$ws.Range("Q13").Value2 = 18.94900735972
$ws.Range("R13").Value2 = 170.54106623748
$ws.Range("S13").Value2 = 0.005583392208333766
$ws.Range("T13").Value2 = 0.005583392208333767

# Row 14: Resolving-Mac -> ECs
$ws.Range("A14").Value2 = "Resolving-Mac"
$ws.Range("B14").Value2 = "B2m"
$ws.Range("C14").Value2 = "Gm11127"
$ws.Range("D14").Value2 = "ECs"
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 2852.789998333333
$ws.Range("H14").Value2 = 8558.369995
$ws.Range("I14").Value2 = 0.4712768877470153
$ws.Range("J14").Value2 = 0.4712768877470154
$ws.Range("K14").Value2 = 2
$ws.Range("L14").Value2 = 0.6666666666666666
$ws.Range("M14").Value2 = 0.4271833333333333
$ws.Range("N14").Value2 = 1.28155
$ws.Range("O14").Value2 = 0.7619379997514806
$ws.Range("P14").Value2 = 0.7619379997514806
$ws.Range("Q14").Value2 = 1218.664340788028
$ws.Range("R14").Value2 = 10967.97906709225
$ws.Range("S14").Value2 = 0.3590837691790639
$ws.Range("T14").Value2 = 0.3590837691790639

# Row 15: Resolving-Mac -> FAPs
$ws.Range("A15").Value2 = "Resolving-Mac"
$ws.Range("B15").Value2 = "B2m"
$ws.Range("C15").Value2 = "Gm11127"
$ws.Range("D15").Value2 = "FAPs"
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 2852.789998333333
$ws.Range("H15").Value2 = 8558.369995
$ws.Range("I15").Value2 = 0.4712768877470153
$ws.Range("J15").Value2 = 0.4712768877470154
$ws.Range("K15").Value2 = 1
$ws.Range("L15").Value2 = 0.3333333333333333
$ws.Range("M15").Value2 = 0.03181266666666666
$ws.Range("N15").Value2 = 0.095438
$ws.Range("O15").Value2 = 0.05674210044109227
$ws.Range("P15").Value2 = 0.05674210044109228
$ws.Range("Q15").Value2 = 90.75485728697888
$ws.Range("R15").Value2 = 816.7937155828099
$ws.Range("S15").Value2 = 0.02674124050010651
$ws.Range("T15").Value2 = 0.02674124050010652

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("A16").Value2 = "Resolving-Mac"
$ws.Range("B16").Value2 = "B2m"
$ws.Range("C16").Value2 = "Gm11127"
$ws.Range("D16").Value2 = "MuSCs"
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 2852.789998333333
$ws.Range("H16").Value2 = 8558.369995
$ws.Range("I16").Value2 = 0.4712768877470153
$ws.Range("J16").Value2 = 0.4712768877470154
$ws.Range("K16").Value2 = 1
$ws.Range("L16").Value2 = 0.3333333333333333
$ws.Range("M16").Value2 = 0.02317166666666666
$ws.Range("N16").Value2 = 0.069515
$ws.Range("O16").Value2 = 0.04132973356694953
$ws.Range("P16").Value2 = 0.04132973356694953
$ws.Range("Q16").Value2 = 66.10389891138054
$ws.Range("R16").Value2 = 594.9350902024249
$ws.Range("S16").Value2 = 0.01947774820684532
$ws.Range("T16").Value2 = 0.01947774820684533

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("A17").Value2 = "Resolving-Mac"
$ws.Range("B17").Value2 = "B2m"
$ws.Range("C17").Value2 = "Gm11127"
$ws.Range("D17").Value2 = "Resolving-Mac"
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 2852.789998333333
$ws.Range("H17").Value2 = 8558.369995
$ws.Range("I17").Value2 = 0.4712768877470153
$ws.Range("J17").Value2 = 0.4712768877470154
$ws.Range("K17").Value2 = 1
$ws.Range("L17").Value2 = 0.3333333333333333
$ws.Range("M17").Value2 = 0.078486
$ws.Range("N17").Value2 = 0.235458
$ws.Range("O17").Value2 = 0.1399901662404776
$ws.Range("P17").Value2 = 0.1399901662404776
$ws.Range("Q17").Value2 = 223.90407580919
$ws.Range("R17").Value2 = 2015.13668228271
$ws.Range("S17").Value2 = 0.06597412986099958
$ws.Range("T17").Value2 = 0.0659741298609996
